$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1671377.9
$ws.Cells.Item(17, 10).Value = 1671377.9
$ws.Cells.Item(17, 12).Value = 5014133.699999999
$ws.Cells.Item(17, 14).Value = -5014469.699999999
# Row 19
$ws.Cells.Item(19, 8).Value = 667.4074000000001
$ws.Cells.Item(19, 9).Value = 677.619
$ws.Cells.Item(19, 10).Value = 631.6667
$ws.Cells.Item(19, 11).Value = 677.619
$ws.Cells.Item(19, 12).Value = 631.6667
$ws.Cells.Item(19, 13).Value = -502.619
$ws.Cells.Item(19, 14).Value = -981.6667
# Row 39
$ws.Cells.Item(39, 8).Value = 578.5
$ws.Cells.Item(39, 9).Value = 103.125
$ws.Cells.Item(39, 11).Value = 309.375
$ws.Cells.Item(39, 13).Value = -13.375
# Row 58
$ws.Cells.Item(58, 8).Value = 1197.125
$ws.Cells.Item(58, 9).Value = 115.6
$ws.Cells.Item(58, 10).Value = 2999.6667
$ws.Cells.Item(58, 11).Value = 346.8
$ws.Cells.Item(58, 12).Value = 8999.000100000001
$ws.Cells.Item(58, 13).Value = -196.8
$ws.Cells.Item(58, 14).Value = -9299.000100000001
# Row 70
$ws.Cells.Item(70, 8).Value = 2103.5652
$ws.Cells.Item(70, 9).Value = 932.13336
$ws.Cells.Item(70, 10).Value = 4300
$ws.Cells.Item(70, 11).Value = 2796.40008
$ws.Cells.Item(70, 12).Value = 12900
$ws.Cells.Item(70, 13).Value = -2526.40008
$ws.Cells.Item(70, 14).Value = -13440
# Row 73
$ws.Cells.Item(73, 8).Value = 2103.5652
$ws.Cells.Item(73, 9).Value = 932.13336
$ws.Cells.Item(73, 10).Value = 4300
$ws.Cells.Item(73, 11).Value = 2796.40008
$ws.Cells.Item(73, 12).Value = 12900
$ws.Cells.Item(73, 13).Value = -1860.40008
$ws.Cells.Item(73, 14).Value = -14772
# Row 86
$ws.Cells.Item(86, 8).Value = 64199.125
$ws.Cells.Item(86, 9).Value = 18298.834
$ws.Cells.Item(86, 10).Value = 201900
$ws.Cells.Item(86, 11).Value = 18298.834
$ws.Cells.Item(86, 12).Value = 201900
$ws.Cells.Item(86, 13).Value = -17175.834
$ws.Cells.Item(86, 14).Value = -204146
# Row 87
$ws.Cells.Item(87, 8).Value = 15160.221
$ws.Cells.Item(87, 10).Value = 15160.221
$ws.Cells.Item(87, 12).Value = 15160.221
$ws.Cells.Item(87, 14).Value = -17656.221
# Row 89
$ws.Cells.Item(89, 8).Value = 64199.125
$ws.Cells.Item(89, 9).Value = 18298.834
$ws.Cells.Item(89, 10).Value = 201900
$ws.Cells.Item(89, 11).Value = 91494.17
$ws.Cells.Item(89, 12).Value = 1009500
$ws.Cells.Item(89, 13).Value = -85878.17
$ws.Cells.Item(89, 14).Value = -1020732
# Row 90
$ws.Cells.Item(90, 8).Value = 15160.221
$ws.Cells.Item(90, 10).Value = 15160.221
$ws.Cells.Item(90, 12).Value = 45480.663
$ws.Cells.Item(90, 14).Value = -57960.663
# Row 112
$ws.Cells.Item(112, 8).Value = 1239.5714
$ws.Cells.Item(112, 9).Value = 845
$ws.Cells.Item(112, 10).Value = 1305.3334
$ws.Cells.Item(112, 11).Value = 2535
$ws.Cells.Item(112, 12).Value = 3916.0002
$ws.Cells.Item(112, 13).Value = -1427
$ws.Cells.Item(112, 14).Value = -6132.0002
# Row 115
$ws.Cells.Item(115, 8).Value = 790.2105
$ws.Cells.Item(115, 9).Value = 223.77777
$ws.Cells.Item(115, 11).Value = 671.33331
$ws.Cells.Item(115, 13).Value = 895.66669
# Row 132
$ws.Cells.Item(132, 8).Value = 13390.358
$ws.Cells.Item(132, 9).Value = 14449.905
$ws.Cells.Item(132, 10).Value = 2189.4285
$ws.Cells.Item(132, 11).Value = 43349.715
$ws.Cells.Item(132, 12).Value = 6568.2855
$ws.Cells.Item(132, 13).Value = -40819.715
$ws.Cells.Item(132, 14).Value = -11628.2855
# Row 137
$ws.Cells.Item(137, 8).Value = 1559.7
$ws.Cells.Item(137, 9).Value = 2027.7142
$ws.Cells.Item(137, 10).Value = 1307.6923
$ws.Cells.Item(137, 11).Value = 6083.142599999999
$ws.Cells.Item(137, 12).Value = 3923.0769
$ws.Cells.Item(137, 13).Value = -3533.142599999999
$ws.Cells.Item(137, 14).Value = -9023.0769
# Row 138
$ws.Cells.Item(138, 8).Value = 5372.5776
$ws.Cells.Item(138, 9).Value = 3443.125
$ws.Cells.Item(138, 10).Value = 6074.197
$ws.Cells.Item(138, 11).Value = 10329.375
$ws.Cells.Item(138, 12).Value = 18222.591
$ws.Cells.Item(138, 13).Value = -5189.375
$ws.Cells.Item(138, 14).Value = -28502.591

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Cells.Item(23, 8).Value = 47990.285
$ws.Cells.Item(23, 9).Value = 46672.668
$ws.Cells.Item(23, 10).Value = 48978.5
$ws.Cells.Item(23, 11).Value = 46672.668
$ws.Cells.Item(23, 12).Value = 48978.5
$ws.Cells.Item(23, 14).Value = -49496.5
$ws.Cells.Item(23, 13).Value = -46413.668
# Row 32
$ws.Cells.Item(32, 8).Value = 12679.625
$ws.Cells.Item(32, 9).Value = 9693.25
$ws.Cells.Item(32, 11).Value = 9693.25
$ws.Cells.Item(32, 13).Value = -9406.25
# Row 132
$ws.Cells.Item(132, 8).Value = 3139.037
$ws.Cells.Item(132, 9).Value = 2934.4614
$ws.Cells.Item(132, 10).Value = 3329
$ws.Cells.Item(132, 11).Value = 8803.3842
$ws.Cells.Item(132, 12).Value = 9987
$ws.Cells.Item(132, 13).Value = -6273.3842
$ws.Cells.Item(132, 14).Value = -15047
# Row 134
$ws.Cells.Item(134, 8).Value = 47990
$ws.Cells.Item(134, 10).Value = 57320
$ws.Cells.Item(134, 12).Value = 57320
$ws.Cells.Item(134, 14).Value = -67460

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 255.5
$ws.Cells.Item(22, 9).Value = 170.63637
$ws.Cells.Item(22, 11).Value = 170.63637
$ws.Cells.Item(22, 13).Value = 2.363630000000001
# Row 94
$ws.Cells.Item(94, 8).Value = 642.1429000000001
$ws.Cells.Item(94, 9).Value = 669.7646999999999
$ws.Cells.Item(94, 10).Value = 524.75
$ws.Cells.Item(94, 11).Value = 669.7646999999999
$ws.Cells.Item(94, 12).Value = 524.75
$ws.Cells.Item(94, 13).Value = -218.7646999999999
$ws.Cells.Item(94, 14).Value = -1426.75
# Row 140
$ws.Cells.Item(140, 8).Value = 59796.668
$ws.Cells.Item(140, 10).Value = 59796.668
$ws.Cells.Item(140, 12).Value = 59796.668
$ws.Cells.Item(140, 14).Value = -70156.66800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 895.3333
$ws.Cells.Item(16, 9).Value = 1115.5
$ws.Cells.Item(16, 10).Value = 601.7778
$ws.Cells.Item(16, 11).Value = 1115.5
$ws.Cells.Item(16, 12).Value = 601.7778
$ws.Cells.Item(16, 13).Value = -828.5
$ws.Cells.Item(16, 14).Value = -1175.7778
# Row 31
$ws.Cells.Item(31, 8).Value = 5368.936
$ws.Cells.Item(31, 9).Value = 2196.7778
$ws.Cells.Item(31, 10).Value = 9651.35
$ws.Cells.Item(31, 11).Value = 2196.7778
$ws.Cells.Item(31, 12).Value = 9651.35
$ws.Cells.Item(31, 13).Value = -1901.7778
$ws.Cells.Item(31, 14).Value = -10241.35
# Row 34
$ws.Cells.Item(34, 8).Value = 5368.936
$ws.Cells.Item(34, 9).Value = 2196.7778
$ws.Cells.Item(34, 10).Value = 9651.35
$ws.Cells.Item(34, 11).Value = 2196.7778
$ws.Cells.Item(34, 12).Value = 9651.35
$ws.Cells.Item(34, 13).Value = -1994.7778
$ws.Cells.Item(34, 14).Value = -10055.35
# Row 113
$ws.Cells.Item(113, 8).Value = 895.3333
$ws.Cells.Item(113, 9).Value = 1115.5
$ws.Cells.Item(113, 10).Value = 601.7778
$ws.Cells.Item(113, 11).Value = 1115.5
$ws.Cells.Item(113, 12).Value = 601.7778
$ws.Cells.Item(113, 13).Value = 1054.5
$ws.Cells.Item(113, 14).Value = -4941.7778
# Row 132
$ws.Cells.Item(132, 8).Value = 2723.4783
$ws.Cells.Item(132, 9).Value = 1576.9166
$ws.Cells.Item(132, 10).Value = 3974.2727
$ws.Cells.Item(132, 11).Value = 4730.7498
$ws.Cells.Item(132, 12).Value = 11922.8181
$ws.Cells.Item(132, 13).Value = -2200.7498
$ws.Cells.Item(132, 14).Value = -16982.8181
# Row 138
$ws.Cells.Item(138, 8).Value = 40302.9
$ws.Cells.Item(138, 10).Value = 40302.9
$ws.Cells.Item(138, 12).Value = 40302.9
$ws.Cells.Item(138, 14).Value = -50582.9
# Row 139
$ws.Cells.Item(139, 8).Value = 35590.582
$ws.Cells.Item(139, 10).Value = 35590.582
$ws.Cells.Item(139, 12).Value = 35590.582
$ws.Cells.Item(139, 14).Value = -45870.582
# Row 140
$ws.Cells.Item(140, 8).Value = 72294.625
$ws.Cells.Item(140, 10).Value = 72294.625
$ws.Cells.Item(140, 12).Value = 72294.625
$ws.Cells.Item(140, 14).Value = -82654.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Cells.Item(23, 8).Value = 67.52631
$ws.Cells.Item(23, 9).Value = 40.833332
$ws.Cells.Item(23, 11).Value = 122.499996
$ws.Cells.Item(23, 13).Value = 112.500004
# Row 122
$ws.Cells.Item(122, 8).Value = 990.57574
$ws.Cells.Item(122, 10).Value = 2799.6667
$ws.Cells.Item(122, 12).Value = 25197.0003
$ws.Cells.Item(122, 14).Value = -30097.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Cells.Item(43, 8).Value = 2405.6667
$ws.Cells.Item(43, 9).Value = 2405.6667
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 2405.6667
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -2254.6667
$ws.Cells.Item(43, 14).ClearContents()
# Row 46
$ws.Cells.Item(46, 8).Value = 12578
$ws.Cells.Item(46, 10).Value = 14997.5
$ws.Cells.Item(46, 12).Value = 14997.5
$ws.Cells.Item(46, 14).Value = -15309.5
# Row 57
$ws.Cells.Item(57, 8).Value = 5055
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value = 2423.2415
$ws.Cells.Item(122, 9).Value = 1725.9
$ws.Cells.Item(122, 10).Value = 3972.889
$ws.Cells.Item(122, 11).Value = 5177.700000000001
$ws.Cells.Item(122, 12).Value = 11918.667
$ws.Cells.Item(122, 13).Value = -2727.700000000001
$ws.Cells.Item(122, 14).Value = -16818.667
# Row 123
$ws.Cells.Item(123, 8).Value = 40610.92
$ws.Cells.Item(123, 10).Value = 40610.92
$ws.Cells.Item(123, 12).Value = 40610.92
$ws.Cells.Item(123, 14).Value = -45510.92
# Row 132
$ws.Cells.Item(132, 8).Value = 2687.0952
$ws.Cells.Item(132, 9).Value = 2397.9614
$ws.Cells.Item(132, 10).Value = 3156.9375
$ws.Cells.Item(132, 11).Value = 7193.8842
$ws.Cells.Item(132, 12).Value = 9470.8125
$ws.Cells.Item(132, 13).Value = -4663.8842
$ws.Cells.Item(132, 14).Value = -14530.8125
# Row 135
$ws.Cells.Item(135, 8).Value = 40600
$ws.Cells.Item(135, 10).Value = 40600
$ws.Cells.Item(135, 12).Value = 40600
$ws.Cells.Item(135, 14).Value = -50740

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 127
$ws.Cells.Item(127, 8).Value = 51186.11
$ws.Cells.Item(127, 10).Value = 51186.11
$ws.Cells.Item(127, 12).Value = 51186.11
$ws.Cells.Item(127, 14).Value = -61106.11
# Row 133
$ws.Cells.Item(133, 8).Value = 23666
$ws.Cells.Item(133, 10).Value = 23666
$ws.Cells.Item(133, 12).Value = 23666
$ws.Cells.Item(133, 14).Value = -28726

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 137
$ws.Cells.Item(137, 8).Value = 56050.4
$ws.Cells.Item(137, 10).Value = 56050.4
$ws.Cells.Item(137, 12).Value = 56050.4
$ws.Cells.Item(137, 14).Value = -66250.39999999999

Write-Host "All updates applied."